$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the existing recipients' displayed addresses (rows 2-4 keep
# their original hyperlinks untouched).
$ws.Range("A2").Value = "vahan.hovhannisyan.h@gmail.com"
$ws.Range("A3").Value = "eurotouram@yahoo.com"
$ws.Range("A4").Value = "eurotouram@gmail.com"

# New recipients appended below the existing list (A6 written before A5
# so new shared strings land in the same order as the authored file).
$ws.Range("A6").Value = "lil-3@mail.ru"
$ws.Range("A5").Value = "lilis88@mail.ru "
$ws.Range("A7").Value = "levon.eurotour@gmail.com"

# A5 and A7 become live mailto links (A6 is left as plain text).
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:lilis88@mail.ru")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:levon.eurotour@gmail.com")

# Re-apply the same look the other hyperlinked cells already use.
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)

# A few extra blank, formatted rows under the list (same style as the
# hyperlinked cells above).
$ws.Range("A8:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights picked up by Excel after the edits.
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 14.45
$ws.Rows.Item(4).RowHeight = 14.45
$ws.Rows.Item(5).RowHeight = 14.45
$ws.Rows.Item(6).RowHeight = 14.45
$ws.Rows.Item(7).RowHeight = 14.45

# Leave the selection where the editor last left it.
[void]$ws.Range("D17").Select()
